# Corrected error in plan: "puzzle-board" -> "drop zone"
#
# We replace the word "puzzle-board" (inside the run that reads
# " of the puzzle-board elements in the ") with "drop zone", splitting
# that run into " of the " / "drop zone" / " elements in the " and
# re-seating the "_GoBack" bookmark immediately after the newly typed
# text, exactly like Word does after an in-place edit.

$d = $word.ActiveDocument

$fullText = $d.Content.Text
$oldWord = "puzzle-board"
$newWord = "drop zone"
$idx = $fullText.IndexOf($oldWord)
$len = $oldWord.Length

# Pin down the run boundaries around the target word first (inserting a
# bookmark forces a run split at that exact character, independent of
# shared formatting/rsid) so the later text edit stays isolated to just
# the "puzzle-board" run instead of being re-merged with its neighbours.
$endRange = $d.Range($idx + $len, $idx + $len)
$d.Bookmarks.Add("_TmpSplitEnd", $endRange)

$startRange = $d.Range($idx, $idx)
$d.Bookmarks.Add("_TmpSplitStart", $startRange)

# Replace "puzzle-board" with "drop zone" inside the now-isolated run.
$target = $d.Range($idx, $idx + $len)
$target.Text = $newWord

# Drop the temporary helper bookmarks; the run split they created
# persists even after they're gone.
$d.Bookmarks("_TmpSplitStart").Delete()
$d.Bookmarks("_TmpSplitEnd").Delete()

# Finally, move "_GoBack" to sit right after the freshly typed "drop
# zone", mirroring where Word leaves it after a manual edit.
$newEnd = $idx + $newWord.Length
$goBackRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)
